$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (F column) values for rows 2-7
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 154
$ws1.Range("F3").Value = 49
$ws1.Range("F4").Value = 255
$ws1.Range("F5").Value = 3941
$ws1.Range("F6").Value = 31
$ws1.Range("F7").Value = 440

# Sheet "全部类型" - update "想去人数" (F column) values for rows 2-5 and 8-9
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 154
$ws4.Range("F3").Value = 49
$ws4.Range("F4").Value = 255
$ws4.Range("F5").Value = 3941
$ws4.Range("F8").Value = 31
$ws4.Range("F9").Value = 440
